$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27 and Row 28 have their species-data columns swapped with each other.
# Save the original row 27 values first.
$a27 = $ws.Range("A27").Value2
$b27 = $ws.Range("B27").Value2
$e27 = $ws.Range("E27").Value2
$f27 = $ws.Range("F27").Value2
$g27 = $ws.Range("G27").Value2
$h27 = $ws.Range("H27").Value2
$z27 = $ws.Range("Z27").Value2
$ab27 = $ws.Range("AB27").Value2

$a28 = $ws.Range("A28").Value2
$b28 = $ws.Range("B28").Value2
$e28 = $ws.Range("E28").Value2
$f28 = $ws.Range("F28").Value2
$g28 = $ws.Range("G28").Value2
$h28 = $ws.Range("H28").Value2
$z28 = $ws.Range("Z28").Value2
$ab28 = $ws.Range("AB28").Value2

$ws.Range("A27").Value2 = $a28
$ws.Range("B27").Value2 = $b28
$ws.Range("E27").Value2 = $e28
$ws.Range("F27").Value2 = $f28
$ws.Range("G27").Value2 = $g28
$ws.Range("H27").Value2 = $h28
$ws.Range("Z27").Value2 = $z28
$ws.Range("AB27").Value2 = $ab28

$ws.Range("A28").Value2 = $a27
$ws.Range("B28").Value2 = $b27
$ws.Range("E28").Value2 = $e27
$ws.Range("F28").Value2 = $f27
$ws.Range("G28").Value2 = $g27
$ws.Range("H28").Value2 = $h27
$ws.Range("Z28").Value2 = $z27
$ws.Range("AB28").Value2 = $ab27
